# "fixed some labels in example forms"
#
# The "osm" sheet (list name / name / label) is missing the "label" column
# for a handful of rows (the city/state entries used by the sample
# buildings-with-POIs form). This copies the existing "name" value into the
# missing "label" cell for each of those rows, same as the source edit:
#   C123 <- B123 (Sacramento)
#   C124 <- B124 (Seattle)
#   C125 <- B125 (Bellingham)
#   C127 <- B127 (CA)
#   C128 <- B128 (WA)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("osm")

$rows = @(123, 124, 125, 127, 128)
foreach ($r in $rows) {
    $nameCell = $ws.Cells.Item($r, 2)   # column B = "name"
    $labelCell = $ws.Cells.Item($r, 3)  # column C = "label"
    $labelCell.Value = $nameCell.Value()
}

# Bring the edited sheet/range into focus, mirroring the author's workflow
# (selecting the newly-filled label cells on the osm sheet).
$ws.Activate() | Out-Null
$ws.Range("C123:C128").Select() | Out-Null
